$d = $word.ActiveDocument

# Start from the very last paragraph of the document (the one ending
# "...without having to squash the tabs into a small heading. ").
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

# ---- 1. blank spacer paragraph ------------------------------------------------
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

# ---- 2. "Image Attributions" heading (bold, 14pt) -----------------------------
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.InsertAfter("Image Attributions")
$cur = $d.Paragraphs.Last.Range
$cur.Font.Bold = $true
$cur.Font.BoldBi = $true
$cur.Font.Size = 14
$cur.Font.SizeBi = 14
$cur.Collapse(0)

# ---- 3. blank spacer paragraph -------------------------------------------------
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
# Clear the heading formatting so it does not bleed into the rest of the document
$cur.Font.Bold = $false
$cur.Font.BoldBi = $false
$cur.Font.Size = 12
$cur.Font.SizeBi = 12
$cur.Collapse(0)

# ---- 4. "Home Images:" ----------------------------------------------------------
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.InsertAfter("Home Images:")
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

# ---- 5-8. home image attribution links ------------------------------------------
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.InsertAfter("- https://www.madetobeamomma.com/web-stories/rainbow-painted-pet-rocks/ ")
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.InsertAfter("- https://craftsbyamanda.com/ladybug-painted-rocks/")
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.InsertAfter("- https://www.thebestideasforkids.com/pet-cactus-rocks/")
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.InsertAfter("- https://patch.com/new-jersey/mahwah/freaky-friday-pet-rock-monsters")
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

# ---- 9. blank spacer paragraph ---------------------------------------------------
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

# ---- 10. "Blog Images:" -----------------------------------------------------------
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.InsertAfter("Blog Images:")
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

# ---- 11-13. blog image attribution links -------------------------------------------
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.InsertAfter("- https://nebg.org/2020/04/15/time-talk-pet-rocks/ ")
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.InsertAfter("- https://www.pinterest.com/pin/899664463055746920/")
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.InsertAfter("- https://doodlewash.com/pet-rock-day/")
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

# ---- 14. blank spacer paragraph -----------------------------------------------------
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

# ---- 15. credit line -----------------------------------------------------------------
$cur.InsertParagraphAfter()
$cur = $d.Paragraphs.Last.Range
$cur.InsertAfter("* Some photos included in the site were taken directly by our team *")
$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

# ---- 16. trailing blank paragraph -----------------------------------------------------
$cur.InsertParagraphAfter()

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
